# Increment the ranking table's percentage-like values (columns E and F,
# rows 2-7) by scaling them from fraction (0-1) to percentage (0-100).
# This matches the commit "pontos notáveis - incremento na tabela de
# ranking" where the stored numeric values in E2:F7 are multiplied by 100
# while keeping their existing cell styles (percentage number format) as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = $cell.Value2 * 100
    }
}
